$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add new "E" (scientific-notation) styled, otherwise empty cells in
#    column N on several existing rows (style matches the one already used
#    for B10/B11/B13-B16/B18.../B24/B32, i.e. numFmt "0.00E+00").
# ---------------------------------------------------------------------------
foreach ($r in 12,13,15,16,17,18,26) {
    $ws.Range("N$r").NumberFormat = "0.00E+00"
}

# ---------------------------------------------------------------------------
# 2. Insert the new rows required for the Mualem Van Genuchten parameter
#    table (water / clay / peat columns alongside the existing sand / silt
#    ones). Row numbers below refer to the *current* sheet layout at the
#    moment each Insert() runs.
# ---------------------------------------------------------------------------
$ws.Rows("39:39").Insert()   # new row for alpha_water (before alpha_sand)
$ws.Rows("42:44").Insert()   # new rows for alpha_clay, alpha_peat, n_water
$ws.Rows("47:49").Insert()   # new rows for n_clay, n_peat, residual_wc_water
$ws.Rows("52:53").Insert()   # new rows for residual_wc_clay, residual_wc_peat

# ---------------------------------------------------------------------------
# 3. Populate the alpha_* block (rows 39-44)
# ---------------------------------------------------------------------------
$ws.Range("A39").Value2 = "alpha_water"
$ws.Range("B39").Value2 = 400

$ws.Range("B40").Value2 = 4.0599999999999996

$ws.Range("A42").Value2 = "alpha_clay"
$ws.Range("B42").Value2 = 1.49

$ws.Range("A43").Value2 = "alpha_peat"
$ws.Range("B43").Value2 = 2.31
$ws.Range("D43").Value2 = "from Hydraulic properties of fen peat soils in Poland, Gnatowski 2010"

$ws.Range("A44").Value2 = "n_water"
$ws.Range("B44").Value2 = 2.5

# ---------------------------------------------------------------------------
# 4. Populate the n_* block (rows 45-49)
# ---------------------------------------------------------------------------
$ws.Range("A47").Value2 = "n_clay"
$ws.Range("B47").Value2 = 1.25

$ws.Range("A48").Value2 = "n_peat"
$ws.Range("B48").Value2 = 1.29

$ws.Range("A49").Value2 = "residual_wc_water"
$ws.Range("B49").Value2 = 0

# ---------------------------------------------------------------------------
# 5. Populate the residual_wc_* block (rows 50-53)
# ---------------------------------------------------------------------------
$ws.Range("A52").Value2 = "residual_wc_clay"
$ws.Range("B52").Value2 = 0

$ws.Range("A53").Value2 = "residual_wc_peat"
$ws.Range("B53").Value2 = 0

# ---------------------------------------------------------------------------
# 6. Fix up the view: selection moves to M22, no pinned top-left cell.
# ---------------------------------------------------------------------------
$ws.Range("M22").Select()
